# Applies the "verbesserungen-git" update:
#   - refresh the cached "last modified" date field (22.07.2024 -> 09.09.2024)
#     everywhere it is cached (slide master, every slide layout, notes master)
#   - rename the "KI-Umgang / SBB ChatGPT" module tile to "KI-Umgang / SBB AI Chat"
#   - move the date-range tile from "30.9. - 4.10." to "30.10. - 6.11."

$p = $ppt.ActivePresentation

# NOTE: boolean-ish COM properties (HasTextFrame, HasText, ...) surface here as
# Int64 (-1 / 0) rather than native PowerShell booleans, so they must be used
# in a truthy `if (...)` test rather than compared with `-eq $true`.

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDateField = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) { $isDateField = $true }
        } catch {
        }
        if ($isDateField -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "22.07.2024") {
                $tr.Text = "09.09.2024"
            }
        }
    }
}

# --- Slide master ---
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# --- Every slide layout belonging to the master ---
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# --- Notes master ---
Update-DatePlaceholder $p.NotesMaster.Shapes

# --- Slide text updates ---
$dash = [char]0x2013
$oldRange = "30.9. " + $dash + " 4.10."
$newRange = "30.10. " + $dash + " 6.11."

$s = $p.Slides.Item(1)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        $t = $tr.Text
        if ($t -eq "KI-Umgang / SBB ChatGPT") {
            $tr.Text = "KI-Umgang / SBB AI Chat"
        } elseif ($t -eq $oldRange) {
            $tr.Text = $newRange
        }
    }
}
